# Trade #67 closed at 2026-02-17 21:11:51 - unknown UNKNOWN +0.000%
# Also records a brand-new open trade (#128 / MarketMaking row) that was
# created right after the close.

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's COM layer
# auto-converting date-looking strings (e.g. "2026-02-17") into date serials.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Helper: write an explicit, present-but-empty text cell (matches source
# rows such as G97/L97 which are serialised as an empty string, not a
# truly blank/absent cell). A leading apostrophe forces text entry; the
# apostrophe itself is not stored as part of the value.
function Set-EmptyTextValue($range) {
    $range.Value = "'"
    $range.Style = "Normal"
}

# -----------------------------------------------------------------
# Sheet: Summary
# -----------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.07   # Current Capital
$wsSummary.Range("B4").Value = 0.86      # Total P&L $
$wsSummary.Range("B6").Value = 95        # Total Trades
$wsSummary.Range("B7").Value = 45        # Winning Trades
$wsSummary.Range("B9").Value = 47.37     # Win Rate %

# -----------------------------------------------------------------
# Sheet: Strategy Status (row 5 = MarketMaking)
# -----------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.07     # Capital
$wsStatus.Range("D5").Value = 62         # Trades
$wsStatus.Range("E5").Value = 0.75       # P&L $
$wsStatus.Range("F5").Value = 1.07       # P&L %
$wsStatus.Range("G5").Value = 50         # Win Rate %

# -----------------------------------------------------------------
# Sheet: All Trades
#   Row 96 -> existing open trade gets closed (early_exit)
#   Row 129 -> new trade row (#128) appended
# -----------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("G96").Value = 0.13476
$wsAll.Range("H96").Value = "CLOSED"
$wsAll.Range("I96").Value = 12.2998
$wsAll.Range("J96").Value = 0.01
$wsAll.Range("K96").Value = 101.07
$wsAll.Range("L96").Value = "early_exit"
$wsAll.Range("M96").Value = 0.13

$wsAll.Range("A129").Value = 128
Set-TextValue $wsAll.Range("B129") "2026-02-17"
$wsAll.Range("C129").Value = "21:11:45"
$wsAll.Range("D129").Value = "MarketMaking"
$wsAll.Range("E129").Value = "DOWN"
$wsAll.Range("F129").Value = 0.12
Set-EmptyTextValue $wsAll.Range("G129")
$wsAll.Range("H129").Value = "OPEN"
$wsAll.Range("I129").Value = 0
$wsAll.Range("J129").Value = 0
$wsAll.Range("K129").Value = 101.0546450978375
Set-EmptyTextValue $wsAll.Range("L129")
$wsAll.Range("M129").Value = 0
$wsAll.Range("N129").Value = 0
$wsAll.Range("O129").Value = 0
$wsAll.Range("P129").Value = 0.6
$wsAll.Range("Q129").Value = "Normal spread capture: 19600 bps"

# -----------------------------------------------------------------
# Sheet: MarketMaking
#   Row 63 -> existing open trade gets closed (early_exit)
#   Row 96 -> new trade row (#128) appended
# -----------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Range("G63").Value = 0.13476
$wsMM.Range("H63").Value = "CLOSED"
$wsMM.Range("I63").Value = 12.2998
$wsMM.Range("J63").Value = 0.01
$wsMM.Range("K63").Value = 101.07
$wsMM.Range("P63").Value = "early_exit"
$wsMM.Range("Q63").Value = 0.13

$wsMM.Range("A96").Value = 128
Set-TextValue $wsMM.Range("B96") "2026-02-17"
$wsMM.Range("C96").Value = "21:11:45"
$wsMM.Range("D96").Value = "MarketMaking"
$wsMM.Range("E96").Value = "DOWN"
$wsMM.Range("F96").Value = 0.12
Set-EmptyTextValue $wsMM.Range("G96")
$wsMM.Range("H96").Value = "OPEN"
$wsMM.Range("I96").Value = 0
$wsMM.Range("J96").Value = 0
$wsMM.Range("K96").Value = 101.0546450978375
$wsMM.Range("L96").Value = 0
$wsMM.Range("M96").Value = 0
$wsMM.Range("N96").Value = 0.6
$wsMM.Range("O96").Value = "Normal spread capture: 19600 bps"
Set-EmptyTextValue $wsMM.Range("P96")
$wsMM.Range("Q96").Value = 0

Write-Host "edit.ps1 completed"
